$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Select + copy the concentration/response pair (cols A:B) from the Example sheet
$ws1.Range("A1:B14").Select()
$ws1.Range("A1:B14").Copy()

# New sheet for the normalised inhibition curve, placed right after "Example"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Inhibition"
$ws2.Range("A1").Select()
$ws2.Paste()

# Normalise column A (concentration) -- fixes the inhibition-curve scaling
$normalised = @(3, 0.3, 1, 3, 10, 30, 3, 100, 300, 3, 0.3, 3, 1, 3)
for ($i = 0; $i -lt $normalised.Count; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $normalised[$i]
}

$ws2.Range("E32").Select()
